$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Profile"
$ws.Range("B2").Value = "//button[@ng-click=`"vm.go('triangular.profile')`"]"
$ws.Range("C2").Value = "click"
$ws.Range("A3").Value = "Address"
$ws.Range("B3").Value = "//input[@ng-model=`"vm.userSetting.address.address`"]"
$ws.Range("C3").Value = "input"
$ws.Range("D3").Value = "703 `"B`" Street"
$ws.Range("A4").Value = "City"
$ws.Range("B4").Value = "//input[@ng-model=`"vm.userSetting.address.city`"]"
$ws.Range("C4").Value = "input"
$ws.Range("D4").Value = "Marysville"
$ws.Range("A5").Value = "Country"
$ws.Range("B5").Value = "//input[@ng-model=`"vm.userSetting.address.county`"]"
$ws.Range("C5").Value = "input"
$ws.Range("D5").Value = "USA"
$ws.Range("A6").Value = "State"
$ws.Range("B6").Value = "//input[@ng-model=`"vm.userSetting.address.state`"]"
$ws.Range("C6").Value = "input"
$ws.Range("D6").Value = "CA "
$ws.Range("A7").Value = "zipcode"
$ws.Range("B7").Value = "//input[@ng-model=`"vm.userSetting.address.zipcode`"]"
$ws.Range("C7").Value = "input"
$ws.Range("D7").Value = 95901
$ws.Range("A8").Value = "phone1"
$ws.Range("B8").Value = "//input[@ng-model=`"vm.userSetting.phone1`"]"
$ws.Range("C8").Value = "input"
$ws.Range("D8").Value = "530-741-4211"
$ws.Range("A9").Value = "phone2"
$ws.Range("B9").Value = "//input[@ng-model=`"vm.userSetting.phone2`"]"
$ws.Range("C9").Value = "input"
$ws.Range("D9").Value = "530-741-4211"
$ws.Range("A10").Value = "email"
$ws.Range("B10").Value = "//input[@ng-model=`"vm.userSetting.email`"]"
$ws.Range("C10").Value = "input"
$ws.Range("D10").Value = "a@gmail.com"
$ws.Range("A11").Value = "Update Settings"
$ws.Range("B11").Value = "(//button[@class=`"md-raised md-primary margin-left-0 md-button ng-scope md-cs-content-theme-theme md-ink-ripple`"])[1]"
$ws.Range("C11").Value = "click"
$ws.Range("A12").Value = "Password"
$ws.Range("B12").Value = "/html/body/div/div/md-content/div[2]/div/md-tabs/md-tabs-wrapper/md-tabs-canvas/md-pagination-wrapper/md-tab-item[2]"
$ws.Range("C12").Value = "click"
$ws.Range("A13").Value = "Current Password"
$ws.Range("B13").Value = "//*[@id=`"old-password`"]"
$ws.Range("C13").Value = "input"
$ws.Range("D13").Value = "Password0!"
$ws.Range("A14").Value = "New Password"
$ws.Range("B14").Value = "//*[@id=`"password`"]"
$ws.Range("C14").Value = "input"
$ws.Range("D14").Value = "Password0!"
$ws.Range("A15").Value = "Confirm Password"
$ws.Range("B15").Value = "//*[@id=`"confirm`"]"
$ws.Range("C15").Value = "input"
$ws.Range("D15").Value = "Password0!"
$ws.Range("A16").Value = "Update Settings"
$ws.Range("B16").Value = "(//button[@class=`"md-raised md-primary margin-left-0 md-button ng-scope md-cs-content-theme-theme md-ink-ripple`"])[2]"
$ws.Range("C16").Value = "click"
$ws.Range("B17").Value = "/html/body/div[3]/md-dialog/md-dialog-actions/button[2]"
$ws.Range("C17").Value = "click"
$ws.Range("A18").Value = "Notification"
$ws.Range("B18").Value = "/html/body/div/div/md-content/div[2]/div/md-tabs/md-tabs-wrapper/md-tabs-canvas/md-pagination-wrapper/md-tab-item[3]"
$ws.Range("C18").Value = "click"
$ws.Range("A19").Value = "Show my location"
$ws.Range("B19").Value = "//md-switch[@aria-label=`"Toggle Show my location`"]"
$ws.Range("C19").Value = "click"
$ws.Range("A20").Value = "show my avatar"
$ws.Range("B20").Value = "//md-switch[@aria-label=`"Toggle Show my avatar`"]"
$ws.Range("C20").Value = "click"
$ws.Range("A21").Value = "send me notifications"
$ws.Range("B21").Value = "//md-switch[@aria-label=`"Toggle Send me notifications`"]"
$ws.Range("C21").Value = "click"
$ws.Range("A22").Value = "show my username"
$ws.Range("B22").Value = "//md-switch[@aria-label=`"Toggle Show my username`"]"
$ws.Range("C22").Value = "click"
$ws.Range("A23").Value = "make my profile public"
$ws.Range("B23").Value = "//md-switch[@aria-label=`"Toggle Make my profile public`"]"
$ws.Range("C23").Value = "click"
$ws.Range("A24").Value = " allow cloud backups"
$ws.Range("B24").Value = "//md-switch[@aria-label=`"Toggle Allow cloud backups`"]"
$ws.Range("C24").Value = "click"
$ws.Range("A25").Value = "Update Settings"
$ws.Range("B25").Value = "(//button[@class=`"md-raised md-primary margin-left-0 md-button ng-scope md-cs-content-theme-theme md-ink-ripple`"])[3]"
$ws.Range("C25").Value = "click"

$ws.Activate()
$ws.Range("B24").Select()

Write-Host "done"
